$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Function signature gains a new keyword argument.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "def send_file(filename, ip, port):", $true, $false, $false, $false, $false,
    $true, 1, $false,
    "def send_file(filename, ip, port, progress_callback=None):", 2) | Out-Null

# ---------------------------------------------------------------------------
# 2. Comment wording tweak.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "# Send the filename first", $true, $false, $false, $false, $false,
    $true, 1, $false,
    "# Send filename first", 2) | Out-Null

# ---------------------------------------------------------------------------
# 3. Replace the "# Then send the file contents" paragraph with the
#    "file_size = os.path.getsize(filename)" statement (multi-run, with the
#    spell-check proofErr wrapper around the dotted call, as in the diff),
#    then add a "sent = 0" paragraph and a blank paragraph after it.
# ---------------------------------------------------------------------------
$pComment = $d.Paragraphs.Item(13)
$xmlFileSize = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t xml:space="preserve">    file_size = </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>os.path.getsize</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>(filename)</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$pComment.Range.InsertXML($xmlFileSize) | Out-Null

$pFileSize = $d.Paragraphs.Item(13)
$pFileSize.Range.InsertParagraphAfter() | Out-Null
$pSent = $d.Paragraphs.Item(14)
$pSent.Range.Text = "    sent = 0"

# Truly empty paragraph (no run at all, matching the diff's bare "<w:p/>"),
# built via InsertXML rather than InsertParagraphAfter (which would leave a
# stray empty <w:r/> behind).
$pSent = $d.Paragraphs.Item(14)
$pSent.Range.InsertParagraphAfter() | Out-Null
$pBlank = $d.Paragraphs.Item(15)
$xmlBlank = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p/></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$pBlank.Range.InsertXML($xmlBlank) | Out-Null

# ---------------------------------------------------------------------------
# 4. After "sock.sendto(data, (ip, port))" add progress-reporting code:
#      sent += len(data)
#      if progress_callback:
#          progress_callback(sent, file_size)
# ---------------------------------------------------------------------------
$pSendData = $d.Paragraphs.Item(21)
$pSendData.Range.InsertParagraphAfter() | Out-Null
$pSentIncr = $d.Paragraphs.Item(22)
$xmlSentIncr = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t xml:space="preserve">            sent += </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>len</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>(data)</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$pSentIncr.Range.InsertXML($xmlSentIncr) | Out-Null

$pSentIncr = $d.Paragraphs.Item(22)
$pSentIncr.Range.InsertParagraphAfter() | Out-Null
$pIfCallback = $d.Paragraphs.Item(23)
$pIfCallback.Range.Text = "            if progress_callback:"

$pIfCallback = $d.Paragraphs.Item(23)
$pIfCallback.Range.InsertParagraphAfter() | Out-Null
$pCallback = $d.Paragraphs.Item(24)
$pCallback.Range.Text = "                progress_callback(sent, file_size)"

# ---------------------------------------------------------------------------
# 5. Drop the "# Send end marker" comment paragraph entirely.
# ---------------------------------------------------------------------------
$pEndMarkerComment = $d.Paragraphs.Item(26)
$pEndMarkerComment.Range.Delete() | Out-Null

# ---------------------------------------------------------------------------
# 6. Append an inline comment to the end-marker send call. Assigned directly
#    (rather than via Find/Replace) so the straight apostrophes around
#    '__END__' survive untouched by smart-quote autocorrect.
# ---------------------------------------------------------------------------
$pEndMarker = $d.Paragraphs.Item(26)
$pEndMarker.Range.Text = "    sock.sendto(b'__END__', (ip, port))  # End of file"

Write-Output "done"
